$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: rename "Position" (C1) to "PositionSize", remove EachCalculateProfit column (D)
$ws.Range("C1").Value = "PositionSize"

# Update B/C values for rows 3-8 (shift data up, insert new row for Stock_Id 3033/180)
$ws.Range("B3").Value = 3033
$ws.Range("C3").Value = 180

$ws.Range("B4").Value = 3035
$ws.Range("C4").Value = 32

$ws.Range("B5").Value = 3141
$ws.Range("C5").Value = 27

$ws.Range("B6").Value = 3189
$ws.Range("C6").Value = 27

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 3504
$ws.Range("C7").Value = 43

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 6411
$ws.Range("C8").Value = 26

# Delete column D (EachCalculateProfit) entirely
$ws.Columns.Item(4).Delete()
